# Update countries & provincias Spain
#
# This script applies the data refresh captured by the commit:
#  - Haiti's stats overtook Gabon's, so the two countries swap rank
#    (row 107 / row 108), with Haiti receiving the newer numbers.
#  - Timor Oriental and Santa Lucia swap rank as well (row 207 / row 208);
#    their totals happen to be identical so only the country names move.
#  - A handful of other countries (Belgica, Kazajistan, Islas Turcas y
#    Caicos, San Martin (Parte Holandesa), Camboya) get refreshed counts.
#  - The "last updated" timestamp banner in A1 is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 05:04"

# --- Belgica (row 35) ---------------------------------------------------
$ws.Range("B35").Value = 114179
$ws.Range("C35").Value = 1376
$ws.Range("D35").Value = 19275
$ws.Range("E35").Value = 84924
$ws.Range("G35").Value = 6
$ws.Range("H35").Value = 9980

# --- Kazajistan (row 39) ------------------------------------------------
$ws.Range("B39").Value = 107775
$ws.Range("C39").Value = 52
$ws.Range("D39").Value = 102736
$ws.Range("E39").Value = 3340

# --- Gabon / Haiti swap rank (rows 107-108) -----------------------------
# Haiti now has more cases than Gabon, so Haiti moves into row 107 with
# its updated numbers, and Gabon (unchanged numbers) drops to row 108.
$ws.Range("A107").Value = "Haiti"
$ws.Range("B107").Value = 8740
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 6688
$ws.Range("E107").Value = 1825
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 227

$ws.Range("A108").Value = "Gabon"
$ws.Range("B108").Value = 8728
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 7934
$ws.Range("E108").Value = 740
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 54

# --- Islas Turcas y Caicos (row 172) ------------------------------------
$ws.Range("B172").Value = 682
$ws.Range("C172").Value = 1
$ws.Range("E172").Value = 57

# --- San Martin (Parte Holandesa) (row 173) -----------------------------
$ws.Range("B173").Value = 644
$ws.Range("C173").Value = 11
$ws.Range("E173").Value = 90

# --- Timor Oriental / Santa Lucia swap rank (rows 207-208) --------------
# Their totals are identical, so only the country names trade places.
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"

# --- Camboya (row 187) ---------------------------------------------------
$ws.Range("D187").Value = 275
$ws.Range("E187").Value = 1
